$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2880.6667
$ws.Range("J17").Value = 2880.6667
$ws.Range("L17").Value = 8642.000100000001
$ws.Range("N17").Value = -8978.000100000001

$ws.Range("H111").Value = 6875
$ws.Range("J111").Value = 6500
$ws.Range("L111").Value = 19500
$ws.Range("N111").Value = -25634

$ws.Range("H127").Value = 2112.1538
$ws.Range("I127").Value = 1086.6666
$ws.Range("J127").Value = 2419.8
$ws.Range("K127").Value = 3259.9998
$ws.Range("L127").Value = 7259.400000000001
$ws.Range("M127").Value = 1700.0002
$ws.Range("N127").Value = -17179.4

$ws.Range("H137").Value = 1404.3513
$ws.Range("I137").Value = 1025.35
$ws.Range("J137").Value = 1850.2354
$ws.Range("K137").Value = 3076.05
$ws.Range("L137").Value = 5550.706200000001
$ws.Range("M137").Value = -526.0499999999997
$ws.Range("N137").Value = -10650.7062

$ws.Range("H138").Value = 551944.75
$ws.Range("J138").Value = 680997.75
$ws.Range("L138").Value = 2042993.25
$ws.Range("N138").Value = -2053273.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1804.2858
$ws.Range("I45").Value = 1773.3334
$ws.Range("K45").Value = 1773.3334
$ws.Range("M45").Value = -1396.3334

$ws.Range("H74").Value = 2117.5
$ws.Range("I74").Value = 1331.6
$ws.Range("J74").Value = 3099.875
$ws.Range("K74").Value = 1331.6
$ws.Range("L74").Value = 3099.875
$ws.Range("M74").Value = -457.5999999999999
$ws.Range("N74").Value = -4847.875

$ws.Range("H77").Value = 2117.5
$ws.Range("I77").Value = 1331.6
$ws.Range("J77").Value = 3099.875
$ws.Range("K77").Value = 6658
$ws.Range("L77").Value = 15499.375
$ws.Range("M77").Value = -2290
$ws.Range("N77").Value = -24235.375

$ws.Range("H122").Value = 2351
$ws.Range("I122").Value = 1838.6666
$ws.Range("K122").Value = 5515.9998
$ws.Range("M122").Value = -3065.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H81").Value = 15980
$ws.Range("J81").Value = 15980
$ws.Range("L81").Value = 15980
$ws.Range("N81").Value = -18102

$ws.Range("H84").Value = 15980
$ws.Range("J84").Value = 15980
$ws.Range("L84").Value = 47940
$ws.Range("N84").Value = -58548

$ws.Range("H94").Value = 19231122
$ws.Range("I94").Value = 41666984
$ws.Range("J94").Value = 382.7143
$ws.Range("K94").Value = 41666984
$ws.Range("L94").Value = 382.7143
$ws.Range("M94").Value = -41666533
$ws.Range("N94").Value = -1284.7143

$ws.Range("H105").Value = 56106364
$ws.Range("I105").Value = 63119416
$ws.Range("J105").Value = 1955
$ws.Range("K105").Value = 63119416
$ws.Range("L105").Value = 1955
$ws.Range("M105").Value = -63117669
$ws.Range("N105").Value = -5449

$ws.Range("H107").Value = 1752.25
$ws.Range("I107").Value = 835
$ws.Range("J107").Value = 2302.6
$ws.Range("K107").Value = 835
$ws.Range("L107").Value = 2302.6
$ws.Range("M107").Value = 1085
$ws.Range("N107").Value = -6142.6

$ws.Range("H110").Value = 19999.5
$ws.Range("J110").Value = 19999.5
$ws.Range("L110").Value = 19999.5
$ws.Range("N110").Value = -28179.5

$ws.Range("H130").Value = 25340
$ws.Range("J130").Value = 25340
$ws.Range("L130").Value = 25340
$ws.Range("N130").Value = -35380

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 50285.07
$ws.Range("I22").Value = 688.5
$ws.Range("J22").Value = 58551.168
$ws.Range("K22").Value = 688.5
$ws.Range("L22").Value = 58551.168
$ws.Range("M22").Value = -338.5
$ws.Range("N22").Value = -59251.168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5557134.5
$ws.Range("J34").Value = 6251767.5
$ws.Range("L34").Value = 18755302.5
$ws.Range("N34").Value = -18755470.5

$ws.Range("H55").Value = 1813
$ws.Range("J55").Value = 1936.5454
$ws.Range("L55").Value = 5809.6362
$ws.Range("N55").Value = -6163.6362

$ws.Range("H60").Value = 1828.2727
$ws.Range("J60").Value = 2288.25
$ws.Range("L60").Value = 6864.75
$ws.Range("N60").Value = -7366.75

$ws.Range("H131").Value = 21277550
$ws.Range("J131").Value = 1161.9412
$ws.Range("L131").Value = 3485.8236
$ws.Range("N131").Value = -13565.8236

$ws.Range("H137").Value = 46884620
$ws.Range("I137").Value = 107145290
$ws.Range("J137").Value = 15214.667
$ws.Range("K137").Value = 321435870
$ws.Range("L137").Value = 45644.001
$ws.Range("M137").Value = -321430770
$ws.Range("N137").Value = -55844.001

$ws.Range("H140").Value = 2780.6875
$ws.Range("I140").Value = 1763.4736
$ws.Range("J140").Value = 3447.138
$ws.Range("K140").Value = 5290.4208
$ws.Range("L140").Value = 10341.414
$ws.Range("M140").Value = -110.4207999999999
$ws.Range("N140").Value = -20701.414

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 10250
$ws.Range("J46").Value = 15000
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15312

$ws.Range("H107").Value = 553.75
$ws.Range("I107").Value = 514.4286
$ws.Range("J107").Value = 584.3333
$ws.Range("K107").Value = 514.4286
$ws.Range("L107").Value = 584.3333
$ws.Range("M107").Value = 1405.5714
$ws.Range("N107").Value = -4424.3333

$ws.Range("H113").Value = 1472
$ws.Range("I113").Value = 1296.375
$ws.Range("K113").Value = 1296.375
$ws.Range("M113").Value = 873.625

$ws.Range("H122").Value = 2733.6155
$ws.Range("I122").Value = 3162
$ws.Range("J122").Value = 1924.4445
$ws.Range("K122").Value = 9486
$ws.Range("L122").Value = 5773.333500000001
$ws.Range("M122").Value = -7036
$ws.Range("N122").Value = -10673.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6385.7144
$ws.Range("I46").Value = 700
$ws.Range("K46").Value = 700
$ws.Range("M46").Value = -512

$ws.Range("H122").Value = 25002108
$ws.Range("I122").Value = 50002092
$ws.Range("J122").Value = 2120.8
$ws.Range("K122").Value = 150006276
$ws.Range("L122").Value = 6362.400000000001
$ws.Range("M122").Value = -150003826
$ws.Range("N122").Value = -11262.4

$ws.Range("H132").Value = 2653.6785
$ws.Range("J132").Value = 2819.6428
$ws.Range("L132").Value = 8458.928400000001
$ws.Range("N132").Value = -13518.9284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 15000
$ws.Range("J33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15500

$ws.Range("H36").Value = 15000
$ws.Range("J36").Value = 15000
$ws.Range("L36").Value = 15000
$ws.Range("N36").Value = -15500

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H40").Value = 18361.334
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H42").Value = 3000
$ws.Range("J42").Value = 3000
$ws.Range("L42").Value = 3000
$ws.Range("N42").Value = -3756

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H136").Value = 979.6
$ws.Range("I136").Value = 869.875
$ws.Range("J136").Value = 1272.2
$ws.Range("K136").Value = 2609.625
$ws.Range("L136").Value = 3816.6
$ws.Range("M136").Value = -59.625
$ws.Range("N136").Value = -8916.6

$ws.Range("H140").Value = 28188.363
$ws.Range("I140").Value = 10000
$ws.Range("J140").Value = 30007.2
$ws.Range("K140").Value = 10000
$ws.Range("L140").Value = 30007.2
$ws.Range("M140").Value = -4820
$ws.Range("N140").Value = -40367.2
